$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph "Dia 18/10: 8min (1 dia)" right before the
#    "Documentação" paragraph, matching the Arial/24 half-pt (12pt), justified,
#    1.5-line-spacing formatting used throughout this list.
# ---------------------------------------------------------------------------
$docIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Documenta*") {
        $docIdx = $i
        break
    }
}

$target = $d.Paragraphs.Item($docIdx)
$target.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($docIdx)
$newPara.Range.Text = "Dia 18/10: 8min (1 dia)"
$newPara.Range.Font.Name = "Arial"
$newPara.Range.Font.Size = 12
$newPara.Format.LineSpacingRule = 5
$newPara.Format.Alignment = 3

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> up one paragraph: off the run holding
#    "Dia 26/09: 1hr (1 dia)" and onto the run holding "Dia 25/09: 36min
#    (1 dia)" immediately before it - pagination recalculated a line earlier.
# ---------------------------------------------------------------------------
$idx25 = $null
$idx26 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Dia 25/09*") { $idx25 = $i }
    if ($t -like "Dia 26/09*") { $idx26 = $i }
}

$p25 = $d.Paragraphs.Item($idx25)
$xml25 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3B7C0B69" w14:textId="2B52D3F7" w:rsidR="00746B69" w:rsidRDefault="00746B69" w:rsidP="00661576"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Dia 25/09: 36min (1 dia)</w:t></w:r></w:p>'
[void]$p25.Range.InsertXML($xml25)

# re-locate "Dia 26/09" (indices shift as InsertXML can re-seat ranges)
$idx26 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Dia 26/09*") { $idx26 = $i; break }
}
$p26 = $d.Paragraphs.Item($idx26)
$xml26 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="339C6AFC" w14:textId="3DC01E1E" w:rsidR="008F4649" w:rsidRDefault="008F4649" w:rsidP="00661576"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Dia 26/09: 1hr (1 dia)</w:t></w:r></w:p>'
[void]$p26.Range.InsertXML($xml26)

Write-Host "edit complete"
